# BIO-iTC-CommentsMatrix.xlsx update
#
# Adds four new selectable "DOCUMENT" values (Eye, Face, Finger, Vein) to the
# Instructions sheet's lookup table (column H) so the individual biometric
# toolboxes can be picked in Review!G1, and documents them in the
# instructions paragraph (Instructions!B3).

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsReview       = $wb.Worksheets.Item("Review")

# --- Populate the new DOCUMENT list entries -------------------------------
# Order matters: each *new* unique string value becomes a new shared-string
# table entry the first time it is written, in write order. Writing Face
# before Eye (even though Face lands in H13 and Eye in H12) reproduces the
# target shared-string ordering (Face, Eye, Finger, Vein).
$wsInstructions.Range("H13").Value = "Face"
$wsInstructions.Range("H12").Value = "Eye"
$wsInstructions.Range("H14").Value = "Finger"
$wsInstructions.Range("H15").Value = "Vein"

# --- Grow the "Table1" DOCUMENT list table to cover the new rows ----------
$documentTable = $wsInstructions.ListObjects.Item(1)
$documentTable.Resize($wsInstructions.Range("H7:H15"))

# --- Update the instructions paragraph with the new toolbox bullets -------
$instructionsText = "The specific document the comments are for should be selected in G1. The fields are:`n- cPP - for the PP-Module`n- CFG - for the PP-Configuration`n- SD - for the Supporting Document`n- TB - for the PAD Toolbox overview`n- Eye - for the Eye Toolbox`n- Face - for the Face Toolbox`n- Finger - for the Fingerprint Toolbox`n- Vein - for the Vein Toolbox"
$wsInstructions.Range("B3").Value = $instructionsText

# The extra lines make the wrapped cell taller; match the authored height.
$wsInstructions.Rows.Item(3).RowHeight = 141.75

# --- Point the Review!G1 dropdown validation at the expanded range --------
$validation = $wsReview.Range("G1").Validation
$validation.Modify(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    "=Instructions!`$H`$8:`$H`$15"
)
